{"js": "// Apply the four textual updates described by the diff:\n//   1) Arabic spelled-out amount:  \u0648\u0627\u062d\u062f \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u062a\u0633\u0639\u0645\u0626\u0629 \u0648\u0639\u0634\u0631\u0648\u0646 \u0623\u0644\u0641  ->  \u0623\u0631\u0628\u0639\u0629 \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u0645\u0626\u062a\u0627\u0646 \u0648\u0633\u062a\u0648\u0646 \u0623\u0644\u0641\n//   2) Numeric amount:             41 920 000,00                          ->  44 260 000,00\n//   3) Month name:                 \u062f\u064a\u0633\u0645\u0628\u0631                                  ->  \u062c\u0627\u0646\u0641\u064a\n//   4) Year (both occurrences):    2020                                    ->  2021\n//\n// \"\u062f\u064a\u0633\u0645\u0628\u0631\" also occurs earlier in the document inside an unrelated historical\n// reference (\"...18 \u062f\u064a\u0633\u0645\u0628\u0631 1991...\") that must stay untouched, so every\n// candidate hit is checked against its own enclosing paragraph's full text\n// before being rewritten (the historical-reference paragraph mentions \"1991\"\n// and never mentions the grant amount \"41 920 000,00\").\n\nconst body = context.document.body;\n\nasync function replaceExact(find, replace, { skipIfParagraphContains } = {}) {\n  const hits = body.search(find, { matchCase: true, matchWholeWord: false });\n  hits.load(\"items\");\n  await context.sync();\n\n  for (const hit of hits.items) {\n    if (skipIfParagraphContains) {\n      const paras = hit.paragraphs;\n      paras.load(\"items/text\");\n      await context.sync();\n      const paraText = paras.items.length ? paras.items[0].text : \"\";\n      if (paraText.includes(skipIfParagraphContains)) continue;\n    }\n    hit.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Unique, unambiguous strings in the document body \u2014 a plain search finds\n// exactly the run(s) the diff touches.\nawait replaceExact(\"\u0648\u0627\u062d\u062f \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u062a\u0633\u0639\u0645\u0626\u0629 \u0648\u0639\u0634\u0631\u0648\u0646 \u0623\u0644\u0641\", \"\u0623\u0631\u0628\u0639\u0629 \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u0645\u0626\u062a\u0627\u0646 \u0648\u0633\u062a\u0648\u0646 \u0623\u0644\u0641\");\nawait replaceExact(\"41 920 000,00\", \"44 260 000,00\");\nawait replaceExact(\"2020\", \"2021\");\n\n// \"\u062f\u064a\u0633\u0645\u0628\u0631\" additionally appears inside an unrelated historical-law sentence\n// (\"...18 \u062f\u064a\u0633\u0645\u0628\u0631 1991...\") \u2014 skip any hit whose paragraph mentions \"1991\".\nawait replaceExact(\"\u062f\u064a\u0633\u0645\u0628\u0631\", \"\u062c\u0627\u0646\u0641\u064a\", { skipIfParagraphContains: \"1991\" });\n", "ps1": "# Apply the four textual updates described by the diff:\n#   1) Arabic spelled-out amount:  \u0648\u0627\u062d\u062f \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u062a\u0633\u0639\u0645\u0626\u0629 \u0648\u0639\u0634\u0631\u0648\u0646 \u0623\u0644\u0641  ->  \u0623\u0631\u0628\u0639\u0629 \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u0645\u0626\u062a\u0627\u0646 \u0648\u0633\u062a\u0648\u0646 \u0623\u0644\u0641\n#   2) Numeric amount:             41 920 000,00                          ->  44 260 000,00\n#   3) Month name:                 \u062f\u064a\u0633\u0645\u0628\u0631                                  ->  \u062c\u0627\u0646\u0641\u064a\n#   4) Year (both occurrences):    2020                                    ->  2021\n#\n# \"\u062f\u064a\u0633\u0645\u0628\u0631\" also appears earlier in the document inside an unrelated historical\n# reference (\"...18 \u062f\u064a\u0633\u0645\u0628\u0631 1991...\") that must stay untouched, so replacements\n# are scoped to the individual paragraph(s) that actually carry the grant\n# amount/date being revised, found dynamically by their unique text, rather\n# than a document-wide Find/Replace.\n\n$d = $word.ActiveDocument\n\n$amountWords = \"\u0648\u0627\u062d\u062f \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u062a\u0633\u0639\u0645\u0626\u0629 \u0648\u0639\u0634\u0631\u0648\u0646 \u0623\u0644\u0641\"\n$amountWordsNew = \"\u0623\u0631\u0628\u0639\u0629 \u0648\u0623\u0631\u0628\u0639\u0648\u0646 \u0645\u0644\u064a\u0648\u0646 \u0648\u0645\u0626\u062a\u0627\u0646 \u0648\u0633\u062a\u0648\u0646 \u0623\u0644\u0641\"\n$amountDigits = \"41 920 000,00\"\n$amountDigitsNew = \"44 260 000,00\"\n$monthOld = \"\u062f\u064a\u0633\u0645\u0628\u0631\"\n$monthNew = \"\u062c\u0627\u0646\u0641\u064a\"\n$yearOld = \"2020\"\n$yearNew = \"2021\"\n\nfunction Replace-InRange($range, $findText, $replaceText) {\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs($i)\n    $t = $para.Range.Text\n\n    # The \"Article 1\" paragraph carries the spelled-out amount, the numeric\n    # amount, the month name AND the first year \u2014 all four live in this one\n    # paragraph.\n    if ($t.IndexOf($amountWords) -ge 0) {\n        Replace-InRange $para.Range $amountWords $amountWordsNew\n        Replace-InRange $para.Range $amountDigits $amountDigitsNew\n        Replace-InRange $para.Range $monthOld $monthNew\n        Replace-InRange $para.Range $yearOld $yearNew\n    }\n    # The \"Article 2\" paragraph repeats the budget year on its own.\n    elseif ($t.IndexOf($yearOld) -ge 0) {\n        Replace-InRange $para.Range $yearOld $yearNew\n    }\n}\n"}
